$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Rename the Pearson logo picture (footer, both Primary and First-page
# headers/footers) from "image2.png" to "image1.png".
$ftrPrimary = $sec.Footers.Item(1)
$ftrPrimary.Range.InlineShapes.Item(1).Name = "image1.png"

$ftrFirstPage = $sec.Footers.Item(2)
$ftrFirstPage.Range.InlineShapes.Item(1).Name = "image1.png"

# Rename the BTEC logo picture (header, both Primary and First-page
# headers/footers) from "image1.jpg" to "image2.jpg".
$hdrPrimary = $sec.Headers.Item(1)
$hdrPrimary.Range.InlineShapes.Item(1).Name = "image2.jpg"

$hdrFirstPage = $sec.Headers.Item(2)
$hdrFirstPage.Range.InlineShapes.Item(1).Name = "image2.jpg"
